# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to match refreshed output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 4..12 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 241
$ws1.Range("F5").Value = 2848
$ws1.Range("F6").Value = 1980
$ws1.Range("F7").Value = 379
$ws1.Range("F8").Value = 133
$ws1.Range("F9").Value = 1045
$ws1.Range("F11").Value = 281
$ws1.Range("F12").Value = 43

# Sheet "全部类型" (all types) - rows 4..13 in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 241
$ws4.Range("F5").Value = 2848
$ws4.Range("F6").Value = 1980
$ws4.Range("F7").Value = 379
$ws4.Range("F9").Value = 133
$ws4.Range("F10").Value = 1045
$ws4.Range("F12").Value = 281
$ws4.Range("F13").Value = 43
